# Auto-generated edit script: updates crypto price/volume data in Sheet1
# per the commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.266.48"
$ws.Range("E2").Value = "  -0.65%  "
$ws.Range("D3").Value = "1.676.74"
$ws.Range("E3").Value = "  -1.31%  "
$ws.Range("E4").Value = "  -0.76%  "
$ws.Range("D5").Value = "'212.06"
$ws.Range("E5").Value = "  -3.20%  "
$ws.Range("D6").Value = "'0.5276"
$ws.Range("E6").Value = "  -3.64%  "
$ws.Range("E7").Value = "  -0.73%  "
$ws.Range("D8").Value = "'0.2660"
$ws.Range("E8").Value = "  -3.04%  "
$ws.Range("D9").Value = "'0.06296"
$ws.Range("E9").Value = "  -2.34%  "
$ws.Range("D10").Value = "'21.40"
$ws.Range("E10").Value = "  -2.96%  "
$ws.Range("D11").Value = "'0.07559"
$ws.Range("E11").Value = "  -1.88%  "
$ws.Range("D12").Value = "1.671.27"
$ws.Range("E12").Value = "  -1.25%  "
$ws.Range("D13").Value = "'4.474"
$ws.Range("E13").Value = "  -1.65%  "
$ws.Range("D14").Value = "'0.5639"
$ws.Range("E14").Value = "  -3.37%  "
$ws.Range("D15").Value = "'67.21"
$ws.Range("E15").Value = "  +2.19%  "
$ws.Range("D16").Value = "'0.000008038"
$ws.Range("E16").Value = "  -4.31%  "
$ws.Range("D17").Value = "26.027.82"
$ws.Range("E17").Value = "  -1.69%  "
$ws.Range("E18").Value = "  -0.71%  "
$ws.Range("D19").Value = "'4.835"
$ws.Range("E19").Value = "  -2.30%  "
$ws.Range("D20").Value = "'188.51"
$ws.Range("E20").Value = "  -1.44%  "
$ws.Range("D21").Value = "'10.45"
$ws.Range("E21").Value = "  -4.94%  "
$ws.Range("D22").Value = "'6.202"
$ws.Range("E22").Value = "  -0.87%  "
$ws.Range("E23").Value = "  -0.77%  "
$ws.Range("D24").Value = "'150.04"
$ws.Range("E24").Value = "  +0.65%  "
$ws.Range("D25").Value = "'0.1257"
$ws.Range("E25").Value = "  -5.11%  "
$ws.Range("D26").Value = "'7.598"
$ws.Range("E26").Value = "  -3.84%  "
$ws.Range("D27").Value = "'16.03"
$ws.Range("E27").Value = "  +1.39%  "
$ws.Range("D28").Value = "'0.06234"
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").Value = "'1.362"
$ws.Range("E29").Value = "  -1.25%  "
$ws.Range("D30").Value = "'1.286"
$ws.Range("E30").Value = "  -3.50%  "
$ws.Range("D31").Value = "'3.512"
$ws.Range("E31").Value = "  -2.84%  "
$ws.Range("D32").Value = "'3.448"
$ws.Range("E32").Value = "  -4.33%  "
$ws.Range("D33").Value = "'1.637"
$ws.Range("E33").Value = "  -3.16%  "
$ws.Range("D34").Value = "'1.005"
$ws.Range("E34").Value = "  -3.47%  "
$ws.Range("D35").Value = "'0.6078"
$ws.Range("E35").Value = "  -1.61%  "
$ws.Range("D36").Value = "'2.407"
$ws.Range("D37").Value = "'2.734"
$ws.Range("E37").Value = "  -1.36%  "
$ws.Range("D38").Value = "'0.01624"
$ws.Range("E38").Value = "  -1.09%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "1.107.72"
$ws.Range("E39").Value = "  -0.94%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'6.109"
$ws.Range("E40").Value = "  -0.60%  "
$ws.Range("D41").Value = "'0.8720"
$ws.Range("E41").Value = "  -0.84%  "
$ws.Range("E42").Value = "  -0.97%  "
$ws.Range("D43").Value = "'99.95"
$ws.Range("E43").Value = "  -1.37%  "
$ws.Range("D44").Value = "1.826.30"
$ws.Range("E44").Value = "  -1.37%  "
$ws.Range("E45").Value = "  -1.05%  "
$ws.Range("D46").Value = "'56.21"
$ws.Range("E46").Value = "  -2.31%  "
$ws.Range("E47").Value = "  -0.27%  "
$ws.Range("D48").Value = "'8.024"
$ws.Range("E48").Value = "  -2.66%  "
$ws.Range("D49").Value = "'0.05237"
$ws.Range("E49").Value = "  -0.97%  "
$ws.Range("D50").Value = "'0.4257"
$ws.Range("E50").Value = "  -1.06%  "
$ws.Range("D51").Value = "'5.987"
$ws.Range("E51").Value = "  -2.74%  "
